$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values to reflect a new, successful (PASSED) transaction run
$ws.Range("G2").Value = "AAACT231840NRHJDP"
$ws.Range("F2").Value = "PASSED"
$ws.Range("H2").Value = "3 jul. 2023, 10:54:17"

# Update the active selection shown in the sheet view
$ws.Range("H4:H5").Select()
